$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted before the old row 47, shifting the
# former rows 47-49 down to 48-50 (data for those rows is unchanged).
$ws.Rows.Item(47).Insert()

$ws.Range("A47").Value = 8
$ws.Range("B47").Value = "Terminal La Palmera de La Serena"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 45124
$ws.Range("D47").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 100112026
$ws.Range("G47").Value = "Haba"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 400
$ws.Range("K47").Value = 12500
$ws.Range("L47").Value = 13000
$ws.Range("M47").Value = 12750
$ws.Range("N47").Value = "`$/saco 25 kilos"
$ws.Range("O47").Value = "Provincia de Limarí"
$ws.Range("P47").Value = 510
$ws.Range("Q47").Value = 25
$ws.Range("R47").Value = "Hortaliza"
